$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace retailer "Masum Electronics" details with "Achol Telecom" details
$ws.Range("C2").Value = "Achol Telecom"
$ws.Range("D2").Value = "Jonail"
$ws.Range("E2").Value = "Arifur Rahman"
$ws.Range("I2").Value = "Arifur Rahman"
$ws.Range("J2").Value = 1912021212
$ws.Range("N2").Value = "Jonail, Baraigram,Natore."
$ws.Range("P2").Value = 1912021212
$ws.Range("T2").Value = 1912021212

# Row 3: clear retailer name ("SR Telecom" removed)
$ws.Range("C3").Value = ""

# Update the active selection to R14
$ws.Range("R14").Select()
